# Edit summary (per the authoritative OOXML diff):
#   1. Slide 16's table (graphic-frame shape 3) switches from the
#      (now-unused) custom table style {61F6AFF3-7D9F-4CC6-B381-EFCCD22B3087}
#      to the built-in table style {B0CD9495-8E73-47A6-81E8-0D5404CABAD3}.
#   2. The deck's theme colour scheme changes from the "Integral" palette
#      to the standard "Office" palette (dk1/lt1 are unchanged black/white;
#      the other ten slots move to the Office Theme colours).

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 16 -------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $tableShape = $shp
    }
}
$table = $tableShape.Table
$table.ApplyStyle("{B0CD9495-8E73-47A6-81E8-0D5404CABAD3}")

# --- 2. Swap the theme colour scheme over to the Office Theme palette ------
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

$colors.Item(1).RGB  = 0         # dk1      000000
$colors.Item(2).RGB  = 16777215  # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388   # dk2      44546A
$colors.Item(4).RGB  = 15132391  # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939  # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501   # accent2  ED7D31
$colors.Item(7).RGB  = 10855845  # accent3  A5A5A5
$colors.Item(8).RGB  = 49407     # accent4  FFC000
$colors.Item(9).RGB  = 12874308  # accent5  4472C4
$colors.Item(10).RGB = 4697456   # accent6  70AD47
$colors.Item(11).RGB = 12673797  # hlink    0563C1
$colors.Item(12).RGB = 7491477   # folHlink 954F72
